# Updates "想去人数" (F column) counters across the four sheets, and
# applies the day's new event list changes on "全部类型" (sheet4):
#   - a new row is inserted for "北京·人气声优 青山渚 专场活动" (2024-10-02)
#   - the now-superseded "北京·伦敦西区音乐剧明星演唱会（摇滚版）" row is removed

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param($ws, $map)
    foreach ($row in $map.Keys) {
        $ws.Range("F$row").Value2 = $map[$row]
    }
}

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
Set-FValues $ws1 @{
    2  = 253
    6  = 555
    8  = 29
    10 = 384
    11 = 353
    12 = 683
    13 = 757
    14 = 1513
    15 = 1513
    16 = 888
    18 = 1352
    19 = 162
    20 = 323
    23 = 102
    24 = 6611
    25 = 4967
    28 = 206
    29 = 200
    32 = 1284
    34 = 251
    35 = 613
    37 = 1337
    38 = 248
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
Set-FValues $ws2 @{
    2  = 13
    4  = 12
    6  = 38
    13 = 6
    18 = 240
}

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
Set-FValues $ws3 @{
    3 = 2457
    4 = 197
    5 = 58
}

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# F-value updates for rows that sit above the inserted row (row 32)
Set-FValues $ws4 @{
    2  = 253
    5  = 13
    6  = 197
    7  = 58
    9  = 555
    11 = 29
    13 = 12
    14 = 384
    15 = 353
    16 = 683
    17 = 757
    18 = 1513
    19 = 1513
    20 = 888
    22 = 1352
    23 = 162
    24 = 323
    26 = 102
    27 = 38
    29 = 6611
    30 = 4967
}

# Insert a brand-new row 32 for "北京·人气声优 青山渚 专场活动"; this pushes
# the former rows 32-40 down to 33-41.
$ws4.Rows.Item(32).Insert()

# Re-create the index/style of column A (it copies formatting+value, then we
# overwrite the value with the correct sequential index).
$ws4.Range("A31").Copy($ws4.Range("A32"))
$ws4.Range("A32").Value2 = 31

# New event details for row 32. Column B holds a date-like string that must
# stay literal text (not get auto-parsed into a date serial number), so we
# force a text number format before assigning it.
$ws4.Range("B32").NumberFormat = "@"
$ws4.Range("B32").Value2 = "2024-10-02"
$ws4.Range("C32").Value2 = "北京·人气声优 青山渚 专场活动"
$ws4.Range("D32").Value2 = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws4.Range("E32").Value2 = "2024.10.02 11:50-10.02 15:40"
$ws4.Range("F32").Value2 = 200
$ws4.Range("G32").Value2 = 458
$ws4.Range("H32").Value2 = "https://show.bilibili.com/platform/detail.html?id=91249"
$ws4.Range("I32").Value2 = "//i2.hdslb.com/bfs/openplatform/202408/xHqpdFa41724641733192.png"

# Fix up the sequential index in column A for the rows that were shifted down
# by the insert (they keep their old values otherwise).
for ($row = 33; $row -le 40; $row++) {
    $ws4.Range("A$row").Value2 = $row - 1
}

# The shifted rows otherwise keep their own previous "想去人数" counter, but a
# few of those events also received their own independent count update
# (matching the same events' updates on the 展览/演出 sheets).
Set-FValues $ws4 @{
    34 = 1284
    36 = 251
    39 = 613
}

# Remove the row that is now displaced to position 41: "北京·伦敦西区音乐剧
# 明星演唱会（摇滚版）" is dropped from this sheet.
$ws4.Rows.Item(41).Delete()

# F-value updates for the rows below the insert/delete pair (their row
# numbers are unaffected because the insertion and deletion cancel out).
Set-FValues $ws4 @{
    42 = 1337
    43 = 248
    49 = 240
}
